$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (TPM-updated) for rows 2-7. Row 1 (header) is unchanged.
# The old "ECs" sending-cluster block (old rows 2-4) is dropped entirely;
# the old "FAPs" (old rows 5-7) and "MuSCs" (old rows 8-10) sending-cluster
# blocks become the new rows 2-4 and 5-7, with updated values.
$data = @(
    @("FAPs", "Angptl1", "Tek", "ECs", 3, 1, 63.34571333333334, 190.03714, 0.9647337009103245, 0.9647337009103245, 3, 1, 53.457377, 160.372131, 0.6217639481372091, 0.6217639481372091, 3386.295678993927, 30476.66111094534, 0.5998366347790248, 0.5998366347790248),
    @("FAPs", "Angptl1", "Tek", "FAPs", 3, 1, 63.34571333333334, 190.03714, 0.9647337009103245, 0.9647337009103245, 3, 1, 31.083557, 93.25067100000001, 0.3615335470438062, 0.3615335470438062, 1969.010091102327, 17721.09081992094, 0.3487835968428081, 0.3487835968428081),
    @("FAPs", "Angptl1", "Tek", "MuSCs", 3, 1, 63.34571333333334, 190.03714, 0.9647337009103245, 0.9647337009103245, 3, 1, 1.436030666666667, 4.308092, 0.01670250481898457, 0.01670250481898457, 90.96638694854224, 818.6974825368801, 0.01611346928849152, 0.01611346928849152),
    @("MuSCs", "Angptl1", "Tek", "ECs", 3, 1, 2.315632666666667, 6.946898, 0.03526629908967548, 0.03526629908967548, 3, 1, 53.457377, 160.372131, 0.6217639481372091, 0.6217639481372091, 123.7876484555153, 1114.088836099638, 0.02192731335818429, 0.02192731335818429),
    @("MuSCs", "Angptl1", "Tek", "FAPs", 3, 1, 2.315632666666667, 6.946898, 0.03526629908967548, 0.03526629908967548, 3, 1, 31.083557, 93.25067100000001, 0.3615335470438062, 0.3615335470438062, 71.97809998539535, 647.802899868558, 0.01274995020099813, 0.01274995020099813),
    @("MuSCs", "Angptl1", "Tek", "MuSCs", 3, 1, 2.315632666666667, 6.946898, 0.03526629908967548, 0.03526629908967548, 3, 1, 1.436030666666667, 4.308092, 0.01670250481898457, 0.01670250481898457, 3.325319522068444, 29.927875698616, 0.000589035530493056, 0.000589035530493056)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# Remove the now-unused trailing rows (old rows 8-10), shrinking the sheet
# from A1:T10 down to A1:T7.
$ws.Rows(8).Resize(3).Delete()
